# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp in A1 / shared string.
# - Update case counters (Casos totales, Nuevos casos, Casos activos,
#   Recuperados, Casos criticos, Muertes hoy, Muertes) for a batch of
#   countries whose figures moved since the last refresh.
# - Re-rank a few countries: "San Martin (Parte Holandesa)" now reports
#   more cases than "Islas Turcas y Caicos"/"Butan"/"Liechtenstein" and
#   moves above them; "Groenlandia" swaps places with "Islas Malvinas"
#   (tied totals, alphabetic/source-order tie-break changed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner (row 1) ------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 17:59"

# --- Straight numeric refreshes (country stays on the same row) -------------
# row => Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$updates = @{
    4   = @(4272500, 24173, 2035233, 2088424, 0, 353, 148843)   # Estados Unidos
    6   = @(1383172, 46150, 883977,  467113,  0, 676, 32082)    # India
    11  = @(343592,  2288,  316169,  18403,   0, 106, 9020)     # Chile
    13  = @(298681,  767,   0,       0,       0, 61,  45738)    # Reino Unido
    17  = @(245864,  274,   198320,  12442,   0, 5,   35102)    # Italia
    21  = @(206096,  136,   190400,  6495,    0, 0,   9201)     # Alemania
    44  = @(52732,   137,   0,       0,       0, 0,   6139)     # Paises Bajos
    71  = @(15130,   49,    9590,    5171,    0, 0,   369)      # Chequia
    78  = @(13248,   555,   5966,    7073,    0, 9,   209)      # Etiopia
    80  = @(10306,   213,   3282,    6950,    0, 4,   74)       # Estado de Palestina
    104 = @(4166,    31,    1374,    2591,    0, 0,   201)      # Grecia
    114 = @(2769,    5,     2103,    655,     0, 0,   11)       # Sri Lanka
    117 = @(2503,    0,     1907,    473,     0, 0,   123)      # Mali
    130 = @(1768,    3,     1297,    405,     0, 0,   66)       # Sierra Leona
    135 = @(1616,    26,    543,     1062,    0, 0,   11)       # Mozambique
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# --- Re-ranked rows (country name also changes on these rows) ---------------
# row => Pais, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$reranked = @{
    187 = @("San Martin (Parte Holandesa)", 93, 9, 63, 15, 0, 0, 15)
    188 = @("Islas Turcas y Caicos",        92, 2, 28, 62, 0, 0, 2)
    189 = @("Butan",                        92, 0, 85, 7,  0, 0, 0)
    190 = @("Liechtenstein",                86, 0, 81, 4,  0, 0, 1)
    210 = @("Groenlandia",                  13, 0, 13, 0,  0, 0, 0)
    211 = @("Islas Malvinas",               13, 0, 13, 0,  0, 0, 0)
}

foreach ($row in $reranked.Keys) {
    $vals = $reranked[$row]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}
